# Update countries & provincias Spain
# Applies the data refresh from "31 de Marzo de 2020 a las 22:50" to "... 23:20"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 23:20"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 185159
$ws.Range("C4").Value = 21371
$ws.Range("D4").Value = 6347
$ws.Range("E4").Value = 175039
$ws.Range("G4").Value = 632
$ws.Range("H4").Value = 3773

# Row 6: España
$ws.Range("B6").Value = 95923
$ws.Range("C6").Value = 7967
$ws.Range("E6").Value = 68200
$ws.Range("G6").Value = 748
$ws.Range("H6").Value = 8464

# Row 8: Alemania
$ws.Range("B8").Value = 71690
$ws.Range("C8").Value = 4805
$ws.Range("D8").Value = 16100
$ws.Range("E8").Value = 54816
$ws.Range("F8").Value = 2675
$ws.Range("G8").Value = 129
$ws.Range("H8").Value = 774

# Row 16: Austria
$ws.Range("B16").Value = 10180
$ws.Range("C16").Value = 562
$ws.Range("E16").Value = 8957

# Row 21: Israel
$ws.Range("F21").Value = 94

# Row 72: Bosnia y Herzegovina
$ws.Range("B72").Value = 420
$ws.Range("C72").Value = 52
$ws.Range("E72").Value = 390
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 13
